$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 101: 115. Distinct Subsequences (Hard, Dynamic Programming)
$ws.Range("A101").Value = "115. Distinct Subsequences"
$ws.Range("B101").Value = "Hard"
$ws.Range("C101").Value = "Dynamic Programming"
$ws.Range("D101").Value = "We can use a top-down solution. The crux is to understand the state transition possibilities: Base case 1: if t is empty, there is only 1 possible s, which is also empty, Base case 2: if t is non-empty, yet s is empty, there can be no solutions, so 0. If s[i] equalts t[j]: check (i+1, j+1) and check (i+1, j). Else if no match, then check (i+1, j). In recursion, these iterations are handled implicitly by the function calls."
$ws.Hyperlinks.Add($ws.Range("E101"), "https://leetcode.com/problems/distinct-subsequences/solutions/37327/easy-to-understand-dp-in-java/ ") | Out-Null

# Row 102: 371. Sum of Two Integers (Medium, Bit Manipulation)
$ws.Range("A102").Value = "371. Sum of Two Integers"
$ws.Range("B102").Value = "Medium"
$ws.Range("C102").Value = "Bit Manipulation"
$ws.Range("D102").Value = "Each digit must be XOR'ed, but the carry is needed in the case that both digits are 1. This is given by (a & b) << 1. We need to perform a ^ b first (XOR), then a & b << 1. We need a temp int for the left shift, as we need to perform it on the old value, not after the XOR."
$ws.Hyperlinks.Add($ws.Range("E102"), "https://leetcode.com/problems/sum-of-two-integers/solutions/84278/a-summary-how-to-use-bit-manipulation-to-solve-problems-easily-and-efficiently/ ") | Out-Null

$ws.Range("E101:E102").Style = "Hyperlink"

# Expand the table (ListObject) to cover the two new rows
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E102"))

# Apply fill colors matching difficulty convention (OLE BGR colors)
$ws.Range("B101").Interior.Color = 255        # Hard -> red FF0000
$ws.Range("B102").Interior.Color = 49407      # Medium -> orange FFC000

# Update sheet view state (selection/active cell)
$ws.Range("D107").Select() | Out-Null
